# chore: adapt column header formatting to respective input file names
#
# - Rename the comparison-table header suffixes from the generic
#   "_old" / "_new" to the concrete format versions "_FV2210" / "_FV2304".
# - Freeze the header row.
# - Turn the sheet's data range into a proper Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1) Rename header row -------------------------------------------------
# "<Name>_old" -> "<Name>_FV2210" (the older / reference format version)
# "<Name>_new" -> "<Name>_FV2304" (the newer format version)
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# --- 2) Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Convert the data range into an Excel Table --------------------------
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false
